$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "'0.67%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'38.61"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'8.54%"
$ws.Range("E3").Style = "Normal"
$ws.Range("B4").Value = "'LEO"
$ws.Range("B4").Style = "Normal"
$ws.Range("C4").Value = "'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("C4").Style = "Normal"
$ws.Range("D4").Value = "'3.726"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'-3.15%"
$ws.Range("E4").Style = "Normal"
$ws.Range("B5").Value = "'HuobiToken"
$ws.Range("B5").Style = "Normal"
$ws.Range("C5").Value = "'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("C5").Style = "Normal"
$ws.Range("D5").Value = "'5.107"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'1.17%"
$ws.Range("E5").Style = "Normal"
$ws.Range("B6").Value = "'Cronos"
$ws.Range("B6").Style = "Normal"
$ws.Range("C6").Value = "'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("C6").Style = "Normal"
$ws.Range("D6").Value = "'0.08108"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'1.28%"
$ws.Range("E6").Style = "Normal"
$ws.Range("B7").Value = "'FTXToken"
$ws.Range("B7").Style = "Normal"
$ws.Range("C7").Value = "'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("C7").Style = "Normal"
$ws.Range("D7").Value = "'1.959"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'5.36%"
$ws.Range("E7").Style = "Normal"
$ws.Range("B8").Value = "'GateToken"
$ws.Range("B8").Style = "Normal"
$ws.Range("C8").Value = "'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("C8").Style = "Normal"
$ws.Range("D8").Value = "'4.198"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'1.48%"
$ws.Range("E8").Style = "Normal"
$ws.Range("B9").Value = "'KuCoinToken"
$ws.Range("B9").Style = "Normal"
$ws.Range("C9").Value = "'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("C9").Style = "Normal"
$ws.Range("D9").Value = "'7.948"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'2.23%"
$ws.Range("E9").Style = "Normal"
$ws.Range("B10").Value = "'MXToken"
$ws.Range("B10").Style = "Normal"
$ws.Range("C10").Value = "'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("C10").Style = "Normal"
$ws.Range("D10").Value = "'0.9304"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'1.09%"
$ws.Range("E10").Style = "Normal"
$ws.Range("B11").Value = "'LiechtensteinCryptoassetsExchange"
$ws.Range("B11").Style = "Normal"
$ws.Range("C11").Value = "'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("C11").Style = "Normal"
$ws.Range("D11").Value = "'0.1410"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'10.47%"
$ws.Range("E11").Style = "Normal"
$ws.Range("B12").Value = "'WazirX"
$ws.Range("B12").Style = "Normal"
$ws.Range("C12").Value = "'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("C12").Style = "Normal"
$ws.Range("D12").Value = "'0.1959"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'2.34%"
$ws.Range("E12").Style = "Normal"
$ws.Range("B13").Value = "'MandalaExchangeToken"
$ws.Range("B13").Style = "Normal"
$ws.Range("C13").Value = "'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("C13").Style = "Normal"
$ws.Range("D13").Value = "'0.09082"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'-0.42%"
$ws.Range("E13").Style = "Normal"
$ws.Range("B14").Value = "'BitrueCoin"
$ws.Range("B14").Style = "Normal"
$ws.Range("C14").Value = "'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("C14").Style = "Normal"
$ws.Range("D14").Value = "'0.03507"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'2.02%"
$ws.Range("E14").Style = "Normal"
$ws.Range("B15").Value = "'BitMartToken"
$ws.Range("B15").Style = "Normal"
$ws.Range("C15").Value = "'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("C15").Style = "Normal"
$ws.Range("D15").Value = "'0.09814"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'-0.47%"
$ws.Range("E15").Style = "Normal"
$ws.Range("B16").Value = "'BitForexToken"
$ws.Range("B16").Style = "Normal"
$ws.Range("C16").Value = "'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("C16").Style = "Normal"
$ws.Range("D16").Value = "'0.001412"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'0.60%"
$ws.Range("E16").Style = "Normal"
$ws.Range("B17").Value = "'TigerCash"
$ws.Range("B17").Style = "Normal"
$ws.Range("C17").Value = "'https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("C17").Style = "Normal"
$ws.Range("D17").Value = "'0.006148"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'-0.74%"
$ws.Range("E17").Style = "Normal"
$ws.Range("E18").Value = "'3.51%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'0.3462"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'1.28%"
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'0.1327"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'0.47%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'4.809"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'-7.54%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.2456"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'6.46%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.04434"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'0.21%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.001224"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'-0.99%"
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'-0.93%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D27").Value = "'0.0001303"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'4.01%"
$ws.Range("E27").Style = "Normal"
$ws.Range("D39").Value = "'0.02081"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'7.36%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.05132"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'-1.91%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.007462"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'-2.12%"
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'0.06%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.1354"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'0.25%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.002134"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'-1.83%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.009271"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'-3.74%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00006244"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Value = "'0.00000000752"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'0.02%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.003033"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Value = "'0.001602"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'-3.50%"
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.00002104"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'0.02%"
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.0002004"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'0.02%"
$ws.Range("E51").Style = "Normal"
